$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings
# (e.g. "4.995", "240.45") are preserved as text, matching the source data
# which uses inline strings for all Price values, not numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.379.35'
$ws.Range("E2").Value = '  -0.34%  '

$ws.Range("D3").Value = '1.846.65'
$ws.Range("E3").Value = '  -0.14%  '

$ws.Range("D4").Value = '0.9988'
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '240.45'
$ws.Range("E5").Value = '  -0.65%  '

$ws.Range("D6").Value = '0.6285'
$ws.Range("E6").Value = '  +0.15%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").Value = '0.07535'
$ws.Range("E8").Value = '  -0.13%  '

$ws.Range("E9").Value = '  -0.75%  '

$ws.Range("D10").Value = '24.51'
$ws.Range("E10").Value = '  +0.75%  '

$ws.Range("D11").Value = '0.07725'
$ws.Range("E11").Value = '  +0.36%  '

$ws.Range("D12").Value = '1.847.65'
$ws.Range("E12").Value = '  -3.37%  '

$ws.Range("D13").Value = '4.995'

$ws.Range("E14").Value = '  -0.35%  '

$ws.Range("D15").Value = '0.000009989'
$ws.Range("E15").Value = '  +2.25%  '

$ws.Range("D16").Value = '82.86'
$ws.Range("E16").Value = '  -1.17%  '

$ws.Range("D17").Value = '6.157'
$ws.Range("E17").Value = '  -0.91%  '

$ws.Range("D18").Value = '29.418.47'
$ws.Range("E18").Value = '  -0.54%  '

$ws.Range("D19").Value = '227.98'
$ws.Range("E19").Value = '  -2.75%  '

$ws.Range("E20").Value = '  -0.45%  '

$ws.Range("D21").Value = '0.9999'
$ws.Range("E21").Value = '  +0.05%  '

$ws.Range("D22").Value = '7.538'
$ws.Range("E22").Value = '  -1.25%  '

$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.07%  '

$ws.Range("D24").Value = '157.08'
$ws.Range("E24").Value = '  +1.09%  '

$ws.Range("D25").Value = '0.1395'
$ws.Range("E25").Value = '  +0.14%  '

$ws.Range("D26").Value = '8.357'
$ws.Range("E26").Value = '  -0.87%  '

$ws.Range("D27").Value = '17.64'
$ws.Range("E27").Value = '  -0.37%  '

$ws.Range("D28").Value = '1.467'
$ws.Range("E28").Value = '  -0.75%  '

$ws.Range("E29").Value = '  -2.85%  '

$ws.Range("D30").Value = '1.256'
$ws.Range("E30").Value = '  -0.68%  '

$ws.Range("E31").Value = '  +0.66%  '

$ws.Range("D32").Value = '4.017'
$ws.Range("E32").Value = '  -0.15%  '

$ws.Range("D33").Value = '1.843'
$ws.Range("E33").Value = '  -2.15%  '

$ws.Range("E34").Value = '  -1.27%  '

$ws.Range("E35").Value = '  -0.59%  '

$ws.Range("D36").Value = '2.587'
$ws.Range("E36").Value = '  +0.00%  '

$ws.Range("D37").Value = '1.253.45'
$ws.Range("E37").Value = '  +1.47%  '

$ws.Range("D38").Value = '0.01814'
$ws.Range("E38").Value = '  +2.17%  '

$ws.Range("D39").Value = '2.780'
$ws.Range("E39").Value = '  -0.47%  '

$ws.Range("D40").Value = '0.9139'
$ws.Range("E40").Value = '  +0.64%  '

$ws.Range("D41").Value = '6.212'
$ws.Range("E41").Value = '  +1.33%  '

$ws.Range("E42").Value = '  +0.08%  '

$ws.Range("D43").Value = '2.016.32'
$ws.Range("E43").Value = '  -3.34%  '

$ws.Range("D44").Value = '101.18'
$ws.Range("E44").Value = '  -0.74%  '

$ws.Range("E45").Value = '  -1.67%  '

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '7.040'
$ws.Range("E46").Value = '  -3.62%  '

$ws.Range("B47").Value = 'TheSandbox'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D47").Value = '0.4022'
$ws.Range("E47").Value = '  -0.15%  '

$ws.Range("D48").Value = '9.097'
$ws.Range("E48").Value = '  -0.67%  '

$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '1.689'
$ws.Range("E49").Value = '  -1.26%  '

$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '0.1126'
$ws.Range("E50").Value = '  +0.94%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.05735'
$ws.Range("E51").Value = '  -0.29%  '
